# Updated cryptos list values (price + 1h volume change) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.906.92"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.628.80"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.16"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.74"
$ws.Range("E8").Value = "  +10.68%  "
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.17"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.630.49"
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("E14").Value = "  +6.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.19"
$ws.Range("E15").Value = "  +21.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.87"
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.906.46"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.98"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.49"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +4.30%  "
$ws.Range("E23").Value = "  +4.11%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.07"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.59"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.428.23"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  +7.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.553"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.35"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.60"
$ws.Range("E42").Value = "  +9.24%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0500"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +5.71%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.767.90"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "89.50"
$ws.Range("E50").Value = "  +4.24%  "
$ws.Range("E51").Value = "  +4.29%  "
